# Working version of PL model
$wb = $excel.ActiveWorkbook

# Rename the sheets to the PL_ prefixed names
$wb.Worksheets.Item("E1a").Name = "PL_E1a"
$wb.Worksheets.Item("E1b").Name = "PL_E1b"
$wb.Worksheets.Item("E2a").Name = "PL_E2a"

# Move the "active tab" / selection from E1a (PL_E1a) to E2a (PL_E2a)
$wsE2a = $wb.Worksheets.Item("PL_E2a")
$wsE2a.Activate()
